# Update NATMI LR-pair TPM-derived metrics for Ntn1-Unc5c with recalculated values
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.7310083333333334
$ws.Range("H2").Value = 2.193025
$ws.Range("I2").Value = 0.01673731480740535
$ws.Range("J2").Value = 0.01673731480740535
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.003305
$ws.Range("N2").Value = 0.009915
$ws.Range("O2").Value = 0.00777996699693274
$ws.Range("P2").Value = 0.00777996699693274
$ws.Range("Q2").Value = 0.002415982541666667
$ws.Range("R2").Value = 0.021743842875
$ws.Range("S2").Value = 0.0001302157568188873
$ws.Range("T2").Value = 0.0001302157568188873

# Row 3
$ws.Range("G3").Value = 0.7310083333333334
$ws.Range("H3").Value = 2.193025
$ws.Range("I3").Value = 0.01673731480740535
$ws.Range("J3").Value = 0.01673731480740535
$ws.Range("O3").Value = 0.6097422606394873
$ws.Range("P3").Value = 0.6097422606394874
$ws.Range("Q3").Value = 0.1893487025333333
$ws.Range("R3").Value = 1.7041383228
$ws.Range("S3").Value = 0.0102054481677021
$ws.Range("T3").Value = 0.01020544816770211

# Row 4
$ws.Range("G4").Value = 0.7310083333333334
$ws.Range("H4").Value = 2.193025
$ws.Range("I4").Value = 0.01673731480740535
$ws.Range("J4").Value = 0.01673731480740535
$ws.Range("M4").Value = 0.16248
$ws.Range("N4").Value = 0.48744
$ws.Range("O4").Value = 0.3824777723635798
$ws.Range("P4").Value = 0.3824777723635799
$ws.Range("Q4").Value = 0.118774234
$ws.Range("R4").Value = 1.068968106
$ws.Range("S4").Value = 0.006401650882884358
$ws.Range("T4").Value = 0.006401650882884359

# Row 5
$ws.Range("I5").Value = 0.8536212576586365
$ws.Range("J5").Value = 0.8536212576586365
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.003305
$ws.Range("N5").Value = 0.009915
$ws.Range("O5").Value = 0.00777996699693274
$ws.Range("P5").Value = 0.00777996699693274
$ws.Range("Q5").Value = 0.1232177371
$ws.Range("R5").Value = 1.1089596339
$ws.Range("S5").Value = 0.006641145212464411
$ws.Range("T5").Value = 0.006641145212464411

# Row 6
$ws.Range("I6").Value = 0.8536212576586365
$ws.Range("J6").Value = 0.8536212576586365
$ws.Range("O6").Value = 0.6097422606394873
$ws.Range("P6").Value = 0.6097422606394874
$ws.Range("Q6").Value = 9.656989753279998
$ws.Range("R6").Value = 86.91290777951998
$ws.Range("S6").Value = 0.5204889553746993
$ws.Range("T6").Value = 0.5204889553746994

# Row 7
$ws.Range("I7").Value = 0.8536212576586365
$ws.Range("J7").Value = 0.8536212576586365
$ws.Range("M7").Value = 0.16248
$ws.Range("N7").Value = 0.48744
$ws.Range("O7").Value = 0.3824777723635798
$ws.Range("P7").Value = 0.3824777723635799
$ws.Range("Q7").Value = 6.057615105599998
$ws.Range("R7").Value = 54.51853595039999
$ws.Range("S7").Value = 0.3264911570714727
$ws.Range("T7").Value = 0.3264911570714727

# Row 8
$ws.Range("G8").Value = 5.662136666666666
$ws.Range("H8").Value = 16.98641
$ws.Range("I8").Value = 0.129641427533958
$ws.Range("J8").Value = 0.129641427533958
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.003305
$ws.Range("N8").Value = 0.009915
$ws.Range("O8").Value = 0.00777996699693274
$ws.Range("P8").Value = 0.00777996699693274
$ws.Range("Q8").Value = 0.01871336168333333
$ws.Range("R8").Value = 0.16842025515
$ws.Range("S8").Value = 0.001008606027649441
$ws.Range("T8").Value = 0.001008606027649441

# Row 9
$ws.Range("G9").Value = 5.662136666666666
$ws.Range("H9").Value = 16.98641
$ws.Range("I9").Value = 0.129641427533958
$ws.Range("J9").Value = 0.129641427533958
$ws.Range("O9").Value = 0.6097422606394873
$ws.Range("P9").Value = 0.6097422606394874
$ws.Range("Q9").Value = 1.466629287946666
$ws.Range("R9").Value = 13.19966359152
$ws.Range("S9").Value = 0.07904785709708584
$ws.Range("T9").Value = 0.07904785709708585

# Row 10
$ws.Range("G10").Value = 5.662136666666666
$ws.Range("H10").Value = 16.98641
$ws.Range("I10").Value = 0.129641427533958
$ws.Range("J10").Value = 0.129641427533958
$ws.Range("M10").Value = 0.16248
$ws.Range("N10").Value = 0.48744
$ws.Range("O10").Value = 0.3824777723635798
$ws.Range("P10").Value = 0.3824777723635799
$ws.Range("Q10").Value = 0.9199839655999998
$ws.Range("R10").Value = 8.2798556904
$ws.Range("S10").Value = 0.04958496440922273
$ws.Range("T10").Value = 0.04958496440922273

